$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.849.53"
$ws.Range("E2").Value = "  +1.33%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.188.57"
$ws.Range("E3").Value = "  +1.00%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.34"
$ws.Range("E5").Value = "  +2.93%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.52"
$ws.Range("E6").Value = "  +2.45%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.187.42"
$ws.Range("E8").Value = "  +0.96%  "

# Row 9
$ws.Range("E9").Value = "  +1.03%  "

# Row 10
$ws.Range("E10").Value = "  -0.17%  "

# Row 11
$ws.Range("E11").Value = "  -0.75%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.511"
$ws.Range("E12").Value = "  +2.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  +2.09%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.84"
$ws.Range("E14").Value = "  +4.49%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.713.95"
$ws.Range("E15").Value = "  +1.07%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.900.18"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.39"
$ws.Range("E17").Value = "  +3.94%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.198.82"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19
$ws.Range("E19").Value = "  +0.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "509.62"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.29"
$ws.Range("E21").Value = "  +3.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.739"
$ws.Range("E22").Value = "  +2.87%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.15"
$ws.Range("E23").Value = "  -1.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.96"
$ws.Range("E24").Value = "  +2.72%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.73"
$ws.Range("E25").Value = "  +0.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +4.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.98"
$ws.Range("E28").Value = "  +2.37%  "

# Row 29
$ws.Range("E29").Value = "  +4.11%  "

# Row 30
$ws.Range("B30").Value = "Stacks"
$ws.Range("C30").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.88"
$ws.Range("E30").Value = "  +4.47%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.89"
$ws.Range("E31").Value = "  +10.18%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.20"
$ws.Range("E32").Value = "  +1.73%  "

# Row 33
$ws.Range("E33").Value = "  +3.54%  "

# Row 34
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.53"
$ws.Range("E35").Value = "  -0.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.72"
$ws.Range("E36").Value = "  -1.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0900"
$ws.Range("E37").Value = "  +0.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "479.35"
$ws.Range("E38").Value = "  +3.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0416"
$ws.Range("E39").Value = "  -2.12%  "

# Row 40
$ws.Range("E40").Value = "  -3.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.81"
$ws.Range("E41").Value = "  +1.65%  "

# Row 42
$ws.Range("E42").Value = "  +3.36%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.296"
$ws.Range("E43").Value = "  +4.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0648"
$ws.Range("E44").Value = "  +10.33%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.931.48"
$ws.Range("E45").Value = "  -4.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.41"
$ws.Range("E46").Value = "  -0.11%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.38"
$ws.Range("E47").Value = "  -1.14%  "

# Row 48
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("E49").Value = "  +1.27%  "

# Row 50
$ws.Range("E50").Value = "  +3.92%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.61"
$ws.Range("E51").Value = "  +5.23%  "
